$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166610598564148
$ws.Range("B1").Value = 1.623699307441711
$ws.Range("C1").Value = 2.873116731643677
$ws.Range("D1").Value = 3.891113996505737
$ws.Range("E1").Value = 1.56290602684021
